$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 816 (Excel shifts 816..911 down
# to 818..913 automatically, carrying values AND formatting with them).
$ws.Rows("816:817").Insert()

# --- New row 816: "1a plateado" entry for the new weekly date (45142) ---
$ws.Range("A816").Value = 4
$ws.Range("B816").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C816").Value = "Los Lagos"
$ws.Range("D816").Value = 45142
$ws.Range("E816").Value = 10
$ws.Range("F816").Value = "Fruta"
$ws.Range("G816").Value = 100102
$ws.Range("H816").Value = "Cítricos"
$ws.Range("I816").Value = 100102003
$ws.Range("J816").Value = "Limón"
$ws.Range("K816").Value = "Sin especificar"
$ws.Range("L816").Value = "1a plateado"
$ws.Range("M816").Value = 600
$ws.Range("N816").Value = 11000
$ws.Range("O816").Value = 11000
$ws.Range("P816").Value = 11000
$ws.Range("Q816").Value = "$/malla 18 kilos"
$ws.Range("R816").Value = "Provincia de Melipilla"
$ws.Range("S816").Value = 611
$ws.Range("T816").Value = 18

# --- New row 817: "2a plateado" entry for the same new date ---
$ws.Range("A817").Value = 4
$ws.Range("B817").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C817").Value = "Los Lagos"
$ws.Range("D817").Value = 45142
$ws.Range("E817").Value = 10
$ws.Range("F817").Value = "Fruta"
$ws.Range("G817").Value = 100102
$ws.Range("H817").Value = "Cítricos"
$ws.Range("I817").Value = 100102003
$ws.Range("J817").Value = "Limón"
$ws.Range("K817").Value = "Sin especificar"
$ws.Range("L817").Value = "2a plateado"
$ws.Range("M817").Value = 600
$ws.Range("N817").Value = 8000
$ws.Range("O817").Value = 8000
$ws.Range("P817").Value = 8000
$ws.Range("Q817").Value = "$/malla 18 kilos"
$ws.Range("R817").Value = "Provincia de Melipilla"
$ws.Range("S817").Value = 444
$ws.Range("T817").Value = 18
